$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 39: append the next day's data point.
# Force column A to text first so the date-like string "2025/09/30" is
# stored as literal text (matching the existing rows), not auto-converted
# to a date serial number by Excel's input parsing. ClearFormats()
# afterward drops the temporary text-format style so the cell ends up
# with no explicit style, same as its neighbors in row 38.
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "2025/09/30"
$ws.Range("A39").ClearFormats()

$ws.Range("B39").Value = "火"
$ws.Range("C39").Value = 9
$ws.Range("D39").Value = 3
